{"js": "// Apply the LOQ4233.docx content updates:\n//  1. Ativa\u00e7\u00e3o date 2022 -> 2024\n//  2. Remove \"EQD (7), \" from the \"Curso (semestre ideal)\" line\n//  3. Append new PT paragraph text to the \"Programa\" (PT) paragraph\n//  4. Replace the \"Programa\" (EN, italic) paragraph text\n//  5. Replace the \"Norma de recupera\u00e7\u00e3o\" explanatory sentence\n\nconst body = context.document.body;\n\n// --- 1) Ativa\u00e7\u00e3o date -------------------------------------------------\nconst ativacao = body.search(\"Ativa\u00e7\u00e3o: 01/01/2022\", { matchCase: true });\nativacao.load(\"items\");\nawait context.sync();\nif (ativacao.items.length > 0) {\n  ativacao.items[0].insertText(\"Ativa\u00e7\u00e3o: 01/01/2024\", Word.InsertLocation.replace);\n}\n\n// --- 2) Curso (semestre ideal) line -----------------------------------\nconst curso = body.search(\n  \"Curso (semestre ideal): EM (6), EA (5), EB (4), EQD (7), EQN (11)\",\n  { matchCase: true }\n);\ncurso.load(\"items\");\nawait context.sync();\nif (curso.items.length > 0) {\n  curso.items[0].insertText(\n    \"Curso (semestre ideal): EM (6), EA (5), EB (4), EQN (11)\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- 3) \"Programa\" paragraph (Portuguese) -----------------------------\nconst programaPt = body.search(\n  \"1 - A Administra\u00e7\u00e3o das organiza\u00e7\u00f5es - definindo a administra\u00e7\u00e3o2 - O processo administrativo: planejamento, organiza\u00e7\u00e3o, dire\u00e7\u00e3o, controle3 \u2013 Processos de Gest\u00e3o: Marketing, Finan\u00e7as, Gest\u00e3o de Pessoas, Produ\u00e7\u00e3o e Opera\u00e7\u00f5es, Pesquisa e Desenvolvimento, Tecnologia da Informa\u00e7\u00e3o, Log\u00edstica e Meio Ambiente.\",\n  { matchCase: true }\n);\nprogramaPt.load(\"items\");\nawait context.sync();\nif (programaPt.items.length > 0) {\n  programaPt.items[0].insertText(\n    \"1 - A Administra\u00e7\u00e3o das organiza\u00e7\u00f5es - definindo a administra\u00e7\u00e3o2 - O processo administrativo: planejamento, organiza\u00e7\u00e3o, dire\u00e7\u00e3o, controle3 \u2013 Processos de Gest\u00e3o: Marketing, Finan\u00e7as, Gest\u00e3o de Pessoas, Produ\u00e7\u00e3o e Opera\u00e7\u00f5es, Pesquisa e Desenvolvimento, Tecnologia da Informa\u00e7\u00e3o, Log\u00edstica e Meio Ambiente.A disciplina ser\u00e1 ministrada com duas estrat\u00e9gias pedag\u00f3gicas a) aplica\u00e7\u00e3o de diferentes m\u00e9todos ativos para compreender os principais conceitos necess\u00e1rios \u00e0 gest\u00e3o de neg\u00f3cios, e b) aplica\u00e7\u00e3o de conceitos por meio do Programa de Aprendizagem com Extens\u00e3o, por meio do qual o alunos oferecem consultoria a micro e pequenas empresas da regi\u00e3o de Lorena ou de parentes e amigos. Nestas consultorias times de alunos, orientados pelo professor, se debru\u00e7am sobre um pequeno problema de gest\u00e3o da empresa e oferecem solu\u00e7\u00f5es.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- 4) \"Programa\" paragraph (English, italic) ------------------------\nconst programaEn = body.search(\n  \"- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment.\",\n  { matchCase: true }\n);\nprogramaEn.load(\"items\");\nawait context.sync();\nif (programaEn.items.length > 0) {\n  programaEn.items[0].insertText(\n    \"1 - The Administration of organizationsdefining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment. The course will be taught with two pedagogical strategiesa) application of different active methods to understand the main conceptsnecessary for business management, and b) application of concepts through the Extension Learning Program, through which students offer consultancy to micro and small companies in the region of Lorena, or companies of relatives and friends. In these consultanciesteams of students, guided by the professor, work on a small problem ofa small company's management problem and offer solutions.\",\n    Word.InsertLocation.replace\n  );\n}\n\n// --- 5) \"Norma de recupera\u00e7\u00e3o\" explanatory sentence --------------------\nconst recuperacao = body.search(\n  \"NF = (MF + PR)/ 2 , onde NF \u00e9 a m\u00e9dia final da segunda avalia\u00e7\u00e3o, MF \u00e9 a m\u00e9dia final da primeira avalia\u00e7\u00e3o e PR \u00e9 a nota do trabalho de recupera\u00e7\u00e3o\",\n  { matchCase: true }\n);\nrecuperacao.load(\"items\");\nawait context.sync();\nif (recuperacao.items.length > 0) {\n  recuperacao.items[0].insertText(\n    \"Os alunos em recupera\u00e7\u00e3o dever\u00e3o realizar reuni\u00f5es com o professor da disciplina, para orientar na execu\u00e7\u00e3o de um trabalho em formato artigo cient\u00edfico em que se discutam as principais ferramentas de gest\u00e3o e sua aplica\u00e7\u00e3o.\",\n    Word.InsertLocation.replace\n  );\n}\n\nawait context.sync();\n", "ps1": "# Apply the LOQ4233.docx content updates:\n#  1. Ativacao date 2022 -> 2024\n#  2. Remove \"EQD (7), \" from the \"Curso (semestre ideal)\" line\n#  3. Append new PT paragraph text to the \"Programa\" (PT) paragraph\n#  4. Replace the \"Programa\" (EN, italic) paragraph text\n#  5. Replace the \"Norma de recuperacao\" explanatory sentence\n\n$d = $word.ActiveDocument\n\nfunction Replace-UniqueText($doc, $oldText, $newText) {\n    $rng = $doc.Content\n    $rng.Find.ClearFormatting()\n    $rng.Find.Text = $oldText\n    $rng.Find.Forward = $true\n    $rng.Find.Wrap = 1  # wdFindStop\n    $rng.Find.MatchCase = $true\n    $rng.Find.MatchWholeWord = $false\n    $rng.Find.MatchWildcards = $false\n    $found = $rng.Find.Execute()\n    if ($found) {\n        $rng.Text = $newText\n    }\n    return $found\n}\n\n# --- 1) Ativacao date ---\nReplace-UniqueText $d 'Ativa\u00e7\u00e3o: 01/01/2022' 'Ativa\u00e7\u00e3o: 01/01/2024'\n\n# --- 2) Curso (semestre ideal) line ---\nReplace-UniqueText $d 'Curso (semestre ideal): EM (6), EA (5), EB (4), EQD (7), EQN (11)' 'Curso (semestre ideal): EM (6), EA (5), EB (4), EQN (11)'\n\n# --- 3) Programa paragraph (Portuguese) ---\nReplace-UniqueText $d '1 - A Administra\u00e7\u00e3o das organiza\u00e7\u00f5es - definindo a administra\u00e7\u00e3o2 - O processo administrativo: planejamento, organiza\u00e7\u00e3o, dire\u00e7\u00e3o, controle3 \u2013 Processos de Gest\u00e3o: Marketing, Finan\u00e7as, Gest\u00e3o de Pessoas, Produ\u00e7\u00e3o e Opera\u00e7\u00f5es, Pesquisa e Desenvolvimento, Tecnologia da Informa\u00e7\u00e3o, Log\u00edstica e Meio Ambiente.' '1 - A Administra\u00e7\u00e3o das organiza\u00e7\u00f5es - definindo a administra\u00e7\u00e3o2 - O processo administrativo: planejamento, organiza\u00e7\u00e3o, dire\u00e7\u00e3o, controle3 \u2013 Processos de Gest\u00e3o: Marketing, Finan\u00e7as, Gest\u00e3o de Pessoas, Produ\u00e7\u00e3o e Opera\u00e7\u00f5es, Pesquisa e Desenvolvimento, Tecnologia da Informa\u00e7\u00e3o, Log\u00edstica e Meio Ambiente.A disciplina ser\u00e1 ministrada com duas estrat\u00e9gias pedag\u00f3gicas a) aplica\u00e7\u00e3o de diferentes m\u00e9todos ativos para compreender os principais conceitos necess\u00e1rios \u00e0 gest\u00e3o de neg\u00f3cios, e b) aplica\u00e7\u00e3o de conceitos por meio do Programa de Aprendizagem com Extens\u00e3o, por meio do qual o alunos oferecem consultoria a micro e pequenas empresas da regi\u00e3o de Lorena ou de parentes e amigos. Nestas consultorias times de alunos, orientados pelo professor, se debru\u00e7am sobre um pequeno problema de gest\u00e3o da empresa e oferecem solu\u00e7\u00f5es.'\n\n# --- 4) Programa paragraph (English, italic) ---\nReplace-UniqueText $d '- The Administration of organizations - defining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment.' '1 - The Administration of organizationsdefining the administration 2 - The administrative process: planning, organization, direction, control 3 - Management Processes: Marketing, Finance, People Management, Production and Operations, Research and Development, Information Technology, Logistics and Environment. The course will be taught with two pedagogical strategiesa) application of different active methods to understand the main conceptsnecessary for business management, and b) application of concepts through the Extension Learning Program, through which students offer consultancy to micro and small companies in the region of Lorena, or companies of relatives and friends. In these consultanciesteams of students, guided by the professor, work on a small problem ofa small company''s management problem and offer solutions.'\n\n# --- 5) Norma de recuperacao explanatory sentence ---\nReplace-UniqueText $d 'NF = (MF + PR)/ 2 , onde NF \u00e9 a m\u00e9dia final da segunda avalia\u00e7\u00e3o, MF \u00e9 a m\u00e9dia final da primeira avalia\u00e7\u00e3o e PR \u00e9 a nota do trabalho de recupera\u00e7\u00e3o' 'Os alunos em recupera\u00e7\u00e3o dever\u00e3o realizar reuni\u00f5es com o professor da disciplina, para orientar na execu\u00e7\u00e3o de um trabalho em formato artigo cient\u00edfico em que se discutam as principais ferramentas de gest\u00e3o e sua aplica\u00e7\u00e3o.'\n\n"}
